$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Header row (row 1): labels for estimate / time-spent columns change
# (shared-string text for D1 stays the same "Estimación (aproximada)";
#  E1 changes from "Tiempo Empleado (Semanal)" to "Tiempo Empleado (Diario)")
# -----------------------------------------------------------------
$ws.Range("E1").Value = "Tiempo Empleado (Diario)"

# -----------------------------------------------------------------
# Data rows 2-17: tasks reorganised, some renamed, hours re-estimated.
# -----------------------------------------------------------------

# Row 2 - Análisis del proyecto (unchanged text, new "Terminada" fill, E time changed)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Análisis del proyecto"
$ws.Range("C2").Value = "Alta"
$ws.Range("D2").Value = 0.125
$ws.Range("E2").Value = 0.006944444444444444

# Row 3 - Entidad-Relación BBDD
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Entidad-Relación BBDD"
$ws.Range("C3").Value = "Alta"
$ws.Range("D3").Value = 0.041666666666666664
$ws.Range("E3").Value = 0.006944444444444444

# Row 4 - Creación de la BBDD (estimate hours changed)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Creación de la BBDD"
$ws.Range("C4").Value = "Crítica"
$ws.Range("D4").Value = 0.125
$ws.Range("E4").Value = 0.020833333333333332

# Row 5 - Diseño interfaz app (only time spent changes)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Diseño interfaz app"
$ws.Range("C5").Value = "Alta"
$ws.Range("D5").Value = 0.25
$ws.Range("E5").Value = 0.010416666666666666

# Row 6 - Búsqueda de APIS
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Búsqueda de APIS"
$ws.Range("C6").Value = "Media"
$ws.Range("D6").Value = 0.10416666666666667
$ws.Range("E6").Value = 0.006944444444444444

# Row 7 - Implementación API'S (was "Implementación BBDD en código")
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Implementación API'S"
$ws.Range("C7").Value = "Crítica"
$ws.Range("D7").Value = 0.25
$ws.Range("E7").Value = 0.006944444444444444

# Row 8 - Implementación BBDD en código (was "Implementación API'S")
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Implementación BBDD en código"
$ws.Range("C8").Value = "Media"
$ws.Range("D8").Value = 0.625
$ws.Range("E8").Value = 0.020833333333333332

# Row 9 - Implementación de log in (hash, cambiar pass...) (was "Implementación Hibernate")
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Implementación de log in (hash, cambiar pass…)"
$ws.Range("C9").Value = "Media"
$ws.Range("D9").Value = 0.20833333333333334
$ws.Range("E9").Value = 0.010416666666666666

# Row 10 - Lógica operaciones BBDD
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Lógica operaciones BBDD"
$ws.Range("C10").Value = "Alta"
$ws.Range("D10").Value = 0.20833333333333334
$ws.Range("E10").Value = 0.010416666666666666

# Row 11 - Lógica de las API'S
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Lógica de las API'S"
$ws.Range("C11").Value = "Alta"
$ws.Range("D11").Value = 0.41666666666666669
$ws.Range("E11").Value = 0.020833333333333332

# Row 12 - Lógica de el cálculo de estadísticas de juegos
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Lógica de el cálculo de estadísticas de juegos"
$ws.Range("C12").Value = "Alta"
$ws.Range("D12").Value = 0.41666666666666669
$ws.Range("E12").Value = 0.020833333333333332

# Row 13 - Exportación de biblioteca y wishlist (was "Log in en la aplicación (Por decidir)", priority Baja)
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Exportación de biblioteca y wishlist"
$ws.Range("C13").Value = "Baja"
$ws.Range("D13").Value = 0.375
$ws.Range("E13").Value = 0.010416666666666666

# Row 14 - Docketización
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Docketización"
$ws.Range("C14").Value = "Media"
$ws.Range("D14").Value = 0.41666666666666669
$ws.Range("E14").Value = 0.020833333333333332

# Row 15 - Automatización para actualizar precios
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Automatización para actualizar precios"
$ws.Range("C15").Value = "Media"
$ws.Range("D15").Value = 0.625
$ws.Range("E15").Value = 0.020833333333333332

# Row 16 - Documentación
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Documentación"
$ws.Range("C16").Value = "Alta"
$ws.Range("D16").Value = 0.20833333333333334
$ws.Range("E16").Value = 0.006944444444444444

# Row 17 - Pruebas
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Pruebas"
$ws.Range("C17").Value = "Media"
$ws.Range("D17").Value = 0.125
$ws.Range("E17").Value = 0.006944444444444444

# -----------------------------------------------------------------
# Totals row 18 - formula ranges now skip row 6 (D6/E6 excluded) and
# start at row 5 plus rows 7-17.
# -----------------------------------------------------------------
$ws.Range("D18").Formula = "=SUM(D5,D7:D17)"
$ws.Range("E18").Formula = "=SUM(E5,E7:E17)"

# -----------------------------------------------------------------
# New legend row 27 explaining the row colours.
# -----------------------------------------------------------------
$ws.Range("A27").Value = "Leyenda:"
$ws.Range("B27").Value = "En curso"
$ws.Range("C27").Value = "Pausada"
$ws.Range("D27").Value = "Terminada"

# -----------------------------------------------------------------
# Colour coding:
#   - Rows 2,3,4,6 -> "Terminada" (finished) -> blue accent1 fill
#   - Legend B27 "En curso" -> same green fill as the regular task rows
#   - Legend C27 "Pausada"  -> new gold/accent4 fill (not used on any task yet)
#   - Legend D27 "Terminada" -> same blue fill as rows 2,3,4,6
# -----------------------------------------------------------------
foreach ($r in 2,3,4,6) {
    $ws.Range("A$r:C$r").Interior.ThemeColor = 5
    $ws.Range("D$r:E$r").Interior.ThemeColor = 5
}

$ws.Range("B27").Interior.ThemeColor = 10
$ws.Range("C27").Interior.ThemeColor = 8
$ws.Range("D27").Interior.ThemeColor = 5
